{"js": "// Add subsystem abbreviations: Power -> Power (PWR), Propulsion -> Propulsion (PRP),\n// Ground Operation -> Ground Operation (GOP). The trailing \"_GoBack\" bookmark\n// (previously sitting right after \"Command\") moves along with the text so it\n// ends up right after the newly-added \"(GOP)\" abbreviation, before the final period.\n\nconst body = context.document.body;\n\n// 1) \"Power\" -> \"Power (PWR)\"\nlet results = body.search(\"Power\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items\");\nawait context.sync();\nresults.items[0].insertText(\" (PWR)\", Word.InsertLocation.after);\nawait context.sync();\n\n// 2) \"Propulsion\" -> \"Propulsion (PRP)\"\nresults = body.search(\"Propulsion\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items\");\nawait context.sync();\nresults.items[0].insertText(\" (PRP)\", Word.InsertLocation.after);\nawait context.sync();\n\n// 3) \"Ground Operation\" -> \"Ground Operation (GOP)\"\nresults = body.search(\"Ground Operation\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items\");\nawait context.sync();\nresults.items[0].insertText(\" (GOP)\", Word.InsertLocation.after);\nawait context.sync();\n\n// 4) Relocate the \"_GoBack\" bookmark from after \"Command\" to right after \"(GOP)\"\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nresults = body.search(\"Ground Operation (GOP)\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nconst afterGop = results.items[0].getRange(Word.RangeLocation.after);\nafterGop.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Add subsystem abbreviations: Power -> Power (PWR), Propulsion -> Propulsion (PRP),\n# Ground Operation -> Ground Operation (GOP). The trailing \"_GoBack\" bookmark\n# (previously sitting right after \"Command\") moves along with the text so it\n# ends up right after the newly-added \"(GOP)\" abbreviation, before the final period.\n\n$d = $word.ActiveDocument\n\nfunction InsertAfterFind($searchText, $insertText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.Text = $searchText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $true\n    $find.Execute() | Out-Null\n    $range.Collapse($wdCollapseEnd)\n    $range.InsertAfter($insertText)\n}\n\nInsertAfterFind \"Power\" \" (PWR)\"\nInsertAfterFind \"Propulsion\" \" (PRP)\"\nInsertAfterFind \"Ground Operation\" \" (GOP)\"\n\n# Relocate the \"_GoBack\" bookmark from after \"Command\" to right after \"(GOP)\"\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n$range = $d.Content\n$find = $range.Find\n$find.Text = \"Ground Operation (GOP)\"\n$find.MatchCase = $true\n$find.Execute() | Out-Null\n$range.Collapse($wdCollapseEnd)\n$d.Bookmarks.Add(\"_GoBack\", $range)\n"}
